$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the E6 cell value (Unc_Distance) from 1 to 2
$ws.Range("E6").Value = 2

# Update the selection to the full data range A1:E9
$ws.Range("A1:E9").Select()
